$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-28 Saturday" "2024-12-29 Sunday"

Replace-Text "459×2=918" "692×6=4152"
Replace-Text "512×5=2560" "413×8=3304"
Replace-Text "151×3=453" "369×5=1845"
Replace-Text "942×3=2826" "765×6=4590"
Replace-Text "921×5=4605" "427×6=2562"

Replace-Text "218×6=1308" "153×2=306"
Replace-Text "968×4=3872" "297×9=2673"
Replace-Text "385×7=2695" "728×4=2912"
Replace-Text "598×6=3588" "712×7=4984"
Replace-Text "307×6=1842" "885×9=7965"

Replace-Text "293×3=879" "446×3=1338"
Replace-Text "245×7=1715" "652×6=3912"
Replace-Text "162×5=810" "803×7=5621"
Replace-Text "866×6=5196" "151×9=1359"
Replace-Text "255×8=2040" "731×8=5848"

Replace-Text "390×2=780" "209×9=1881"
Replace-Text "294×9=2646" "438×8=3504"
Replace-Text "424×6=2544" "571×9=5139"
Replace-Text "761×6=4566" "320×2=640"
Replace-Text "946×3=2838" "898×5=4490"

Replace-Text "436×4=1744" "578×2=1156"
Replace-Text "179×4=716" "491×7=3437"
Replace-Text "699×6=4194" "234×2=468"
Replace-Text "742×2=1484" "646×9=5814"
Replace-Text "890×7=6230" "818×5=4090"
